$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in new H6 (column N) and T6 (column O) data for rows 2-13
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 14).Value = 45  # Column N -> H6
    $ws.Cells.Item($r, 15).Value = 72  # Column O -> T6
}

# Update the selection to match the new active range
$ws.Range("N2:O13").Select()
